# Applies the "ending position of all balls" update (run #2 timings/trial types)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: Sheet1 -> temp
$ws.Name = "temp"

# Per-row updates: Run Number (A), Run Time (B), Trial Type (F), First Beep (G), Second Beep (H)
$rows = @(
    @{ Row = 2; A = 2; B = "15:55:31.215029"; F = "Kinestic"; G = "15:55:31.221782"; H = "15:55:35.400989" }
    @{ Row = 3; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:55:41.436190"; H = "15:55:46.129023" }
    @{ Row = 4; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:55:52.192796"; H = "15:55:57.238017" }
    @{ Row = 5; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:56:03.318454"; H = "15:56:08.389773" }
    @{ Row = 6; A = 2; B = "15:55:31.215029"; F = "Visual"; G = "15:56:14.453446"; H = "15:56:19.011669" }
    @{ Row = 7; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:56:25.041341"; H = "15:56:30.097322" }
    @{ Row = 8; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:56:36.158616"; H = "15:56:41.231065" }
    @{ Row = 9; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:56:47.278505"; H = "15:56:52.314788" }
    @{ Row = 10; A = 2; B = "15:55:31.215029"; F = "Visual"; G = "15:56:58.399139"; H = "15:57:02.960623" }
    @{ Row = 11; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:57:09.046768"; H = "15:57:14.105278" }
    @{ Row = 12; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:57:20.175255"; H = "15:57:25.218832" }
    @{ Row = 13; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:57:31.364629"; H = "15:57:36.443085" }
    @{ Row = 14; A = 2; B = "15:55:31.215029"; F = "Kinestic"; G = "15:57:42.569074"; H = "15:57:47.186889" }
    @{ Row = 15; A = 2; B = "15:55:31.215029"; F = "Visual"; G = "15:57:53.338331"; H = "15:57:58.000812" }
    @{ Row = 16; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:58:04.156940"; H = "15:58:09.286592" }
    @{ Row = 17; A = 2; B = "15:55:31.215029"; F = "Kinestic"; G = "15:58:15.163829"; H = "15:58:19.537060" }
    @{ Row = 18; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:58:25.418756"; H = "15:58:30.300093" }
    @{ Row = 19; A = 2; B = "15:55:31.215029"; F = "Visual"; G = "15:58:36.433505"; H = "15:58:40.847369" }
    @{ Row = 20; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:58:46.969609"; H = "15:58:52.018465" }
    @{ Row = 21; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:58:58.061267"; H = "15:59:03.097673" }
    @{ Row = 22; A = 2; B = "15:55:31.215029"; F = "Visual"; G = "15:59:09.158147"; H = "15:59:13.708376" }
    @{ Row = 23; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:59:19.802593"; H = "15:59:24.883972" }
    @{ Row = 24; A = 2; B = "15:55:31.215029"; F = "Kinestic"; G = "15:59:31.008496"; H = "15:59:35.685135" }
    @{ Row = 25; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:59:41.746635"; H = "15:59:46.815395" }
    @{ Row = 26; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "15:59:52.899851"; H = "15:59:57.966242" }
    @{ Row = 27; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "16:00:04.066441"; H = "16:00:09.217937" }
    @{ Row = 28; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "16:00:15.121511"; H = "16:00:20.051573" }
    @{ Row = 29; A = 2; B = "15:55:31.215029"; F = "Visual"; G = "16:00:26.122375"; H = "16:00:30.772550" }
    @{ Row = 30; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "16:00:36.881526"; H = "16:00:41.979977" }
    @{ Row = 31; A = 2; B = "15:55:31.215029"; F = "Normal"; G = "16:00:47.937600"; H = "16:00:53.102007" }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
}
